$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.965.41'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.555.32'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'207.13"
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("E6").Value = '  +0.54%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = "'22.10"
$ws.Range("E8").Value = '  +3.80%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  +0.95%  '
$ws.Range("D11").Value = "'0.0857"
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '1.777.51'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '1.555.98'
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("D14").Value = "'3.74"
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("D16").Value = '26.966.33'
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("D17").Value = "'61.69"
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = "'217.95"
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("D23").Value = "'9.23"
$ws.Range("E23").Value = '  +0.42%  '
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("D25").Value = "'154.42"
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("D26").Value = "'6.64"
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").Value = "'14.94"
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = "'0.0468"
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("D33").Value = '1.423.79'
$ws.Range("E33").Value = '  +4.73%  '
$ws.Range("D34").Value = "'3.07"
$ws.Range("E34").Value = '  +4.65%  '
$ws.Range("E35").Value = '  +3.91%  '
$ws.Range("E36").Value = '  +1.88%  '
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("D38").Value = "'0.0165"
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("D39").Value = "'0.520"
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("E41").Value = '  +2.73%  '
$ws.Range("E43").Value = '  +4.50%  '
$ws.Range("D44").Value = "'0.986"
$ws.Range("D45").Value = "'64.41"
$ws.Range("E45").Value = '  +1.61%  '
$ws.Range("E46").Value = '  +1.30%  '
$ws.Range("D47").Value = '1.690.99'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").Value = "'87.83"
$ws.Range("E48").Value = '  +1.95%  '
$ws.Range("E49").Value = '  +2.53%  '
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  +3.55%  '
$ws.Range("E51").Value = '  +0.77%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
